$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O203:O211").Formula = "=C203/B203"
$ws.Range("P203:P211").Formula = "=D203/B203"
$ws.Range("Q204:Q228").Formula = "=D204/C204"
$ws.Range("Q203").Formula = "=D203/C203"

$ws.Range("O212").Formula = "=C212/B212"
$ws.Range("P212").Formula = "=D212/B212"
$ws.Range("R212").Formula = "=E212/B212"
$ws.Range("S212").Formula = "=E212/C212"

$ws.Range("O213:O228").Formula = "=C213/B213"
$ws.Range("P213:P228").Formula = "=D213/B213"
$ws.Range("R213:R228").Formula = "=E213/B213"
$ws.Range("S213:S228").Formula = "=E213/C213"

$ws.Range("O289").Formula = "=C289/B289"
$ws.Range("P289").Formula = "=D289/C289"
$ws.Range("Q289").Formula = "=E289/B289"

$ws.Range("O290").Formula = "=C290/B290"
$ws.Range("P290").Formula = "=D290/C290"
$ws.Range("Q290").Formula = "=E290/B290"
